$d = $word.ActiveDocument

# Locate the last paragraph of the document, i.e. the "GoToAlbumEditPage"
# bullet that currently ends the Controller/Servlets list.
$target = $d.Paragraphs.Last

# --- New top-level bullet: "Filtri" --------------------------------------
$target.Range.InsertParagraphAfter()
$pFiltri = $d.Paragraphs.Last
$pFiltri.Range.Text = "Filtri"
$pFiltri.Style = "Paragrafoelenco"
$pFiltri.Range.ListFormat.ListLevelNumber = 1

# --- New sub-bullet: "LoggedFilter" --------------------------------------
$pFiltri.Range.InsertParagraphAfter()
$pLogged = $d.Paragraphs.Last
$pLogged.Range.Text = "LoggedFilter"
$pLogged.Style = "Paragrafoelenco"
$pLogged.Range.ListFormat.ListLevelNumber = 2

# --- New sub-bullet: "AlreadyLoggedFilter" -------------------------------
$pLogged.Range.InsertParagraphAfter()
$pAlready = $d.Paragraphs.Last
$pAlready.Range.Text = "AlreadyLoggedFilter"
$pAlready.Style = "Paragrafoelenco"
$pAlready.Range.ListFormat.ListLevelNumber = 2

# --- Trailing empty paragraph, indented, no list -------------------------
$pAlready.Range.InsertParagraphAfter()
$pTail = $d.Paragraphs.Last
$pTail.Style = "Normale"
$pTail.Range.ParagraphFormat.LeftIndent = 54
